$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '94.439.11'
$ws.Range('E2').Value = '  +2.38%  '
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.480.50'
$ws.Range('E3').Value = '  +4.88%  '
# Row 4
$ws.Range('E4').Value = '  -0.05%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.32'
$ws.Range('E5').Value = '  +3.76%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '625.92'
$ws.Range('E6').Value = '  +1.87%  '
# Row 7
$ws.Range('E7').Value = '  +6.65%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.394'
$ws.Range('E8').Value = '  +4.62%  '
# Row 9
$ws.Range('E9').Value = '  -0.07%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.00'
$ws.Range('E10').Value = '  +10.40%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.478.83'
$ws.Range('E11').Value = '  +4.91%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.20'
$ws.Range('E12').Value = '  +4.36%  '
# Row 13
$ws.Range('E13').Value = '  +5.70%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.23'
$ws.Range('E14').Value = '  +5.23%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.129.40'
$ws.Range('E15').Value = '  +4.67%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '94.257.68'
$ws.Range('E16').Value = '  +2.46%  '
# Row 17
$ws.Range('E17').Value = '  +4.43%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.37'
$ws.Range('E18').Value = '  +6.72%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.477.90'
$ws.Range('E19').Value = '  +4.71%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.62'
$ws.Range('E20').Value = '  +14.03%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.96'
$ws.Range('E21').Value = '  +5.33%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.503'
$ws.Range('E22').Value = '  +15.66%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '522.19'
$ws.Range('E23').Value = '  +7.94%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.36'
$ws.Range('E24').Value = '  +3.90%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.69'
$ws.Range('E25').Value = '  +10.37%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000185'
$ws.Range('E26').Value = '  +3.55%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '96.27'
$ws.Range('E27').Value = '  +8.31%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.19'
$ws.Range('E28').Value = '  +6.44%  '
# Row 29
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.93'
$ws.Range('E29').Value = '  +11.84%  '
# Row 30
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '11.48'
$ws.Range('E30').Value = '  +4.45%  '
# Row 31
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.140'
$ws.Range('E31').Value = '  +4.84%  '
# Row 32
$ws.Range('E32').Value = '  +0.00%  '
# Row 33
$ws.Range('B33').Value = 'Cronos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.181'
$ws.Range('E33').Value = '  +6.30%  '
# Row 34
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  -0.04%  '
# Row 35
$ws.Range('B35').Value = 'PolygonEcosystemToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.562'
$ws.Range('E35').Value = '  +8.09%  '
# Row 36
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '29.69'
$ws.Range('E36').Value = '  +6.36%  '
# Row 37
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '574.44'
$ws.Range('E37').Value = '  +10.71%  '
# Row 38
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.44'
$ws.Range('E38').Value = '  +8.23%  '
# Row 39
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.56'
$ws.Range('E39').Value = '  +4.53%  '
# Row 40
$ws.Range('B40').Value = 'USDe'
$ws.Range('C40').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.00%  '
# Row 41
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.922'
$ws.Range('E41').Value = '  +6.06%  '
# Row 42
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.149'
$ws.Range('E42').Value = '  +3.28%  '
# Row 43
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0426'
$ws.Range('E43').Value = '  +7.58%  '
# Row 44
$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '23.74'
$ws.Range('E44').Value = '  -1.02%  '
# Row 45
$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.70'
$ws.Range('E45').Value = '  +3.69%  '
# Row 46
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.53'
$ws.Range('E46').Value = '  +4.58%  '
# Row 47
$ws.Range('B47').Value = 'MantraDAO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.54'
$ws.Range('E47').Value = '  -0.15%  '
# Row 48
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.16'
$ws.Range('E48').Value = '  +3.78%  '
# Row 49
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.22'
$ws.Range('E49').Value = '  +4.94%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '53.47'
$ws.Range('E50').Value = '  +3.21%  '
# Row 51
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.06'
$ws.Range('E51').Value = '  +1.87%  '
